$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.775.82"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.292.86"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'301.14"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'99.02"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("D10").Value = "'36.07"
$ws.Range("E10").Value = "  +7.51%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "'17.87"
$ws.Range("E13").Value = "  +9.99%  "
$ws.Range("D14").Value = "'6.84"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "2.648.27"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "2.332.25"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "'0.801"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "42.681.52"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'12.35"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").Value = "'6.19"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'67.84"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "'235.95"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  +13.14%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").Value = "'24.62"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'34.50"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  -5.19%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'4.98"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "'17.47"
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  +3.10%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.101"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "1.998.73"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").Value = "'2.23"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "'10.13"
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.88"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'17.46"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'55.55"
$ws.Range("E48").Value = "  +5.39%  "
$ws.Range("D49").Value = "2.515.37"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").Value = "  -1.18%  "
